$d = $word.ActiveDocument

function Append-ViaFind {
    param($paraIndex, $findText, $replaceText)
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $ok = $r.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Host "FAILED to find/replace in paragraph $paraIndex : [$findText]"
    }
}

# "slime + blocked pieces (jellyfish goal)" -> append " x"
Append-ViaFind 20 "jellyfish goal)" "jellyfish goal) x"

# "slime + blocked pieces + removable obstacle" -> append " x"
Append-ViaFind 21 "slime + blocked pieces + removable obstacle" "slime + blocked pieces + removable obstacle x"

# "time bomb + removable obstacle" -> append " x"
Append-ViaFind 22 "time bomb + removable obstacle" "time bomb + removable obstacle x"

# "locked pieces + removable obstacle" -> append " x" (as two runs " " + "x" in the source edit)
Append-ViaFind 23 "locked pieces + removable obstacle" "locked pieces + removable obstacle x"

# "time bomb + locked pieces + removable obstacle" -> append " x"
Append-ViaFind 24 "locked pieces + removable obstacle" "locked pieces + removable obstacle x"

# "shielded pieces + removable obstacles" -> append " x"
Append-ViaFind 25 "shielded pieces + removable obstacles" "shielded pieces + removable obstacles x"

# "shielded pieces (shielded pieces jellyfish)" -> append " x"
Append-ViaFind 26 "shielded pieces jellyfish)" "shielded pieces jellyfish) x"

# "...create more fields of color" -> "...create more fields of colour"
Append-ViaFind 39 "fields of color" "fields of colour"

# Add seven new paragraphs after the last paragraph ("Swapping removable + obs: ...")
$newParagraphTexts = @(
    "Swapping shielded: used 1/5/10 pieces. Very difficult, and dependent on RNG because of the dependence on bombs in the right positions. Strategic swapping required. On highest shielded count, closest to victory was 3 pieces left",
    "Swapping shielded removable: literally almost impossible without a lot of luck on high piece counts",
    "Swapping blocked removable: probably harder even than shielded pieces, but felt more manageable somehow",
    "Swapping blocked removable time: Time bomb had 10 turns timer. Really difficult as well. Blocked pieces in swap is just a very difficult goal, and removing the bombs gets very difficult in high piece count.",
    "Swapping slime blocked: Turns out blocked pieces are just really difficult to remove without gadgets. Very difficult, and slimes are also more of a threat in swap.",
    "Swapping slime blocked removable: easier than just blocked, maybe because of lower number of blocked pieces.",
    "Swapping time removable: Difficult, but not too much, so was actually pretty fun."
)

foreach ($t in $newParagraphTexts) {
    $newPara = $d.Paragraphs.Add()
    $newPara.Range.Text = $t
}
